$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, copying the format of the existing
# header cell (G1) so it gets the same bold/border/alignment style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add values for the new column in the data rows
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
